# Reorder the "Emission" list values in column E (rows 2-21) of the
# "Lists" worksheet. The set of values is unchanged; only their order
# within the list changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lists")

$newEmissionOrder = @(
    "FERT_ORG",
    "DAPANI",
    "CONVAR",
    "CO2e_sources",
    "RM",
    "CO2e_TRN",
    "RESHID",
    "CO2e_PP",
    "contam_agua",
    "CO2e_AFOLU",
    "Health",
    "CO2e_HFC",
    "CONTUR",
    "CONHAB",
    "CONHICK",
    "turismo_residuos",
    "salud_residuos",
    "CO2e_WASTE",
    "CO2e_DE",
    "CO2e_PIUP"
)

$startRow = 2
for ($i = 0; $i -lt $newEmissionOrder.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 5).Value = $newEmissionOrder[$i]
}
